# "Generate Report for handback" - refresh the handoff/handback timestamps
# for the d2197cfc-... row on the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-25 03:22:56"
$wsZhCn.Range("G4").Value = "2016-01-25 03:25:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-25 03:23:07"
$wsDeDe.Range("G4").Value = "2016-01-25 03:25:32"
